$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 51: tweak the Saturday (2020-09-19) session time + hours ---
$ws.Range("J51").Value = "09:30 - 11:30 13:15 - 12:45 19:30 - 21:30"
$ws.Range("K51").Value = 5.75

# --- Row 52: Sunday 2020-09-20 - "Implementing Flask App" ---
$ws.Range("N52").Value = "Implementing Flask App"
$ws.Range("J52").Value = "09:00-11:00 13:30 - 15:30 20:00 - 22:00"
$ws.Range("K52").Value = 6
$ws.Rows.Item(52).RowHeight = 45

# --- Row 53: Monday 2020-09-21 - "Flask App" ---
$ws.Range("N53").Value = "Flask App"
$ws.Range("J53").Value = "11:00 - 13:00 15:00 - 15:15 21:00 - 23:00"
$ws.Range("K53").Value = 4.25
$ws.Rows.Item(53).RowHeight = 45

# --- Row 54: Tuesday 2020-09-22 - "Flask App" ---
$ws.Range("N54").Value = "Flask App"
$ws.Range("J54").Value = "11:30 - 12:30"
$ws.Range("K54").Value = 1

# --- Row 56: Thursday 2020-09-24 - "Flask App" ---
$ws.Range("N56").Value = "Flask App"
$ws.Range("J56").Value = "12:30 - 14:00"
$ws.Range("K56").Value = 1.5

# --- Row 58: Friday 2020-09-25 - "Flask App + Detection" ---
$ws.Range("N58").Value = "Flask App + Detection"
$ws.Range("J58").Value = "10:30 - 11:30 17:00 - 19:00"
$ws.Range("K58").Value = 1
$ws.Rows.Item(58).RowHeight = 30

# --- Recalculate formulas (weekly sums, totals, averages, chart cache) ---
$excel.Calculate()

# --- Restore window scroll position / selection to match the saved view ---
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 3
$ws.Range("J59").Select() | Out-Null
